$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text formatting so
# numeric-looking strings (e.g. "1.008") are not coerced into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.542.34"
$ws.Range("E2").Value = "  +4.14%  "

$ws.Range("D3").Value = "1.793.05"
$ws.Range("E3").Value = "  +0.74%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").Value = "313.58"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").Value = "0.5319"
$ws.Range("E7").Value = "  +0.86%  "

$ws.Range("D8").Value = "0.3783"
$ws.Range("E8").Value = "  +0.52%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "42.63"
$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.07512"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").Value = "1.115"
$ws.Range("E11").Value = "  +1.76%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "20.96"
$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").Value = "6.160"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").Value = "7.417"
$ws.Range("E15").Value = "  +6.21%  "

$ws.Range("D16").Value = "1.800.34"
$ws.Range("E16").Value = "  +0.94%  "

$ws.Range("D17").Value = "90.81"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "0.00001066"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "0.06447"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "17.18"
$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").Value = "5.915"
$ws.Range("E22").Value = "  +0.31%  "

$ws.Range("D23").Value = "28.673.77"
$ws.Range("E23").Value = "  +4.47%  "

$ws.Range("D24").Value = "11.24"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").Value = "2.130"
$ws.Range("E25").Value = "  +1.55%  "

$ws.Range("D26").Value = "161.35"
$ws.Range("E26").Value = "  +3.48%  "

$ws.Range("D27").Value = "20.46"
$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").Value = "2.398"
$ws.Range("E28").Value = "  +1.17%  "

$ws.Range("D29").Value = "2.020.33"
$ws.Range("E29").Value = "  +1.59%  "

$ws.Range("D30").Value = "123.39"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").Value = "1.120"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").Value = "0.1013"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").Value = "5.684"
$ws.Range("E33").Value = "  +1.48%  "

$ws.Range("D34").Value = "3.659"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").Value = "0.2302"
$ws.Range("E35").Value = "  +12.21%  "

$ws.Range("D36").Value = "0.06596"
$ws.Range("E36").Value = "  +10.12%  "

$ws.Range("D37").Value = "0.02320"
$ws.Range("E37").Value = "  +2.62%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "8.677"
$ws.Range("E38").Value = "  +5.05%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.085"
$ws.Range("E39").Value = "  +3.54%  "

$ws.Range("D40").Value = "11.49"
$ws.Range("E40").Value = "  +1.62%  "

$ws.Range("D41").Value = "0.6316"
$ws.Range("E41").Value = "  +2.85%  "

$ws.Range("D42").Value = "1.200"
$ws.Range("E42").Value = "  +5.72%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "1.404"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").Value = "13.58"
$ws.Range("E45").Value = "  +2.30%  "

$ws.Range("D46").Value = "0.5919"
$ws.Range("E46").Value = "  +2.08%  "

$ws.Range("D47").Value = "3.662"
$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("D48").Value = "124.82"
$ws.Range("E48").Value = "  +2.52%  "

$ws.Range("D49").Value = "1.988"
$ws.Range("E49").Value = "  +4.78%  "

$ws.Range("D50").Value = "1.156"
$ws.Range("E50").Value = "  +3.12%  "

$ws.Range("D51").Value = "0.06917"
$ws.Range("E51").Value = "  +2.68%  "

